$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh re-shuffles which date/volume/price block sits on which
# row (rows 7, 8 and 16 keep their original data). Capture the current
# (pre-edit) values for the columns that move: D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# O (Origen) and P (Precio $/Kg) before overwriting anything. Value2() is
# used so date-formatted cells come back as raw serial numbers instead of
# being coerced into .NET DateTime objects.

$cols = @("D", "J", "K", "L", "M", "O", "P")

$snapshot = @{}
foreach ($r in 2..16) {
    $row = @{}
    foreach ($c in $cols) {
        $cellRef = "$c$r"
        $row[$c] = $ws.Range($cellRef).Value2()
    }
    $snapshot[$r] = $row
}

# Maps destination row -> source row (using the pre-edit snapshot above).
$rowMap = @{
    2  = 14
    3  = 6
    4  = 15
    5  = 10
    6  = 4
    9  = 13
    10 = 9
    11 = 3
    12 = 5
    13 = 11
    14 = 2
    15 = 12
}

foreach ($dst in $rowMap.Keys) {
    $src = $rowMap[$dst]
    $srcRow = $snapshot[$src]
    foreach ($c in $cols) {
        $cellRef = "$c$dst"
        $ws.Range($cellRef).Value = $srcRow[$c]
    }
}
